{"js": "// Replace each three-digit-by-one-digit multiplication equation in the\n// document with its updated counterpart, per the commit diff.\nconst replacements = [\n  [\"687\u00d72=1374\", \"714\u00d72=1428\"],\n  [\"279\u00d77=1953\", \"958\u00d72=1916\"],\n  [\"661\u00d74=2644\", \"797\u00d74=3188\"],\n  [\"280\u00d74=1120\", \"711\u00d75=3555\"],\n  [\"400\u00d79=3600\", \"914\u00d73=2742\"],\n  [\"402\u00d79=3618\", \"979\u00d75=4895\"],\n  [\"497\u00d77=3479\", \"938\u00d77=6566\"],\n  [\"695\u00d77=4865\", \"309\u00d72=618\"],\n  [\"431\u00d77=3017\", \"565\u00d76=3390\"],\n  [\"525\u00d78=4200\", \"666\u00d73=1998\"],\n  [\"142\u00d74=568\", \"287\u00d78=2296\"],\n  [\"961\u00d74=3844\", \"809\u00d76=4854\"],\n  [\"842\u00d75=4210\", \"706\u00d76=4236\"],\n  [\"769\u00d72=1538\", \"685\u00d76=4110\"],\n  [\"178\u00d74=712\", \"908\u00d78=7264\"],\n  [\"694\u00d75=3470\", \"988\u00d76=5928\"],\n  [\"684\u00d77=4788\", \"917\u00d78=7336\"],\n  [\"300\u00d75=1500\", \"268\u00d76=1608\"],\n  [\"698\u00d72=1396\", \"931\u00d72=1862\"],\n  [\"650\u00d74=2600\", \"903\u00d72=1806\"],\n  [\"621\u00d72=1242\", \"667\u00d76=4002\"],\n  [\"632\u00d77=4424\", \"973\u00d77=6811\"],\n  [\"875\u00d74=3500\", \"934\u00d72=1868\"],\n  [\"488\u00d74=1952\", \"409\u00d74=1636\"],\n  [\"600\u00d78=4800\", \"716\u00d74=2864\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication equation in the\n# document with its updated counterpart, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"687\u00d72=1374\", \"714\u00d72=1428\"),\n    @(\"279\u00d77=1953\", \"958\u00d72=1916\"),\n    @(\"661\u00d74=2644\", \"797\u00d74=3188\"),\n    @(\"280\u00d74=1120\", \"711\u00d75=3555\"),\n    @(\"400\u00d79=3600\", \"914\u00d73=2742\"),\n    @(\"402\u00d79=3618\", \"979\u00d75=4895\"),\n    @(\"497\u00d77=3479\", \"938\u00d77=6566\"),\n    @(\"695\u00d77=4865\", \"309\u00d72=618\"),\n    @(\"431\u00d77=3017\", \"565\u00d76=3390\"),\n    @(\"525\u00d78=4200\", \"666\u00d73=1998\"),\n    @(\"142\u00d74=568\", \"287\u00d78=2296\"),\n    @(\"961\u00d74=3844\", \"809\u00d76=4854\"),\n    @(\"842\u00d75=4210\", \"706\u00d76=4236\"),\n    @(\"769\u00d72=1538\", \"685\u00d76=4110\"),\n    @(\"178\u00d74=712\", \"908\u00d78=7264\"),\n    @(\"694\u00d75=3470\", \"988\u00d76=5928\"),\n    @(\"684\u00d77=4788\", \"917\u00d78=7336\"),\n    @(\"300\u00d75=1500\", \"268\u00d76=1608\"),\n    @(\"698\u00d72=1396\", \"931\u00d72=1862\"),\n    @(\"650\u00d74=2600\", \"903\u00d72=1806\"),\n    @(\"621\u00d72=1242\", \"667\u00d76=4002\"),\n    @(\"632\u00d77=4424\", \"973\u00d77=6811\"),\n    @(\"875\u00d74=3500\", \"934\u00d72=1868\"),\n    @(\"488\u00d74=1952\", \"409\u00d74=1636\"),\n    @(\"600\u00d78=4800\", \"716\u00d74=2864\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
